$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; existing rows 45-113 shift down to 46-114,
# carrying their formatting (incl. the date-number style on column D) along.
$ws.Rows.Item(45).Insert()

# Populate the newly-inserted row 45 with a fresh data record.
$ws.Cells.Item(45, 1).Value  = 3
$ws.Cells.Item(45, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(45, 3).Value  = "Coquimbo"
$ws.Cells.Item(45, 4).Value  = 44467
$ws.Cells.Item(45, 5).Value  = 5
$ws.Cells.Item(45, 6).Value  = "Fruta"
$ws.Cells.Item(45, 7).Value  = 100101
$ws.Cells.Item(45, 8).Value  = "Berries"
$ws.Cells.Item(45, 9).Value  = 100112025
$ws.Cells.Item(45, 10).Value = "Frutilla"
$ws.Cells.Item(45, 11).Value = "Sin especificar"
$ws.Cells.Item(45, 12).Value = "Especial"
$ws.Cells.Item(45, 13).Value = 45
$ws.Cells.Item(45, 14).Value = 15000
$ws.Cells.Item(45, 15).Value = 15000
$ws.Cells.Item(45, 16).Value = 15000
$ws.Cells.Item(45, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(45, 19).Value = 2143
$ws.Cells.Item(45, 20).Value = 7
